$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update profile flag values (column C)
$ws.Range("C73").Value = 0
$ws.Range("C110").Value = 1
$ws.Range("C112").Value = 1
$ws.Range("C114").Value = 1

# Update the view: scroll so row 54 is the top-left visible cell, and
# move the active selection to A73.
$ws.Range("A73").Select()
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
